$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column G (the "sum" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

# Add the new header "Save" in H1, copying the style used by the other
# header cells (bold, centered, bordered) from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate H2:H<lastRow> with a binary flag derived from the "sum" (G)
# column: 1 when the performance score exceeds 8, otherwise 0.
for ($r = 2; $r -le $lastRow; $r++) {
    $gVal = $ws.Cells.Item($r, 7).Value2
    if ($gVal -gt 8) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}

$excel.CutCopyMode = 0
